$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto values per the diff (prices & 1h volume %, plus two row swaps).
# Numeric-looking Price (column D) values must stay TEXT cells (as in the source data),
# so we force them with a leading apostrophe and then restore the default 'Normal' style
# (the apostrophe trick otherwise stamps a quotePrefix style onto the cell).
$ws.Range('D2').Value = '62.528.60'
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').Value = '3.019.82'
$ws.Range('E3').Value = '  +2.38%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = "'596.52"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.40%  '
$ws.Range('D6').Value = "'149.55"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.76%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.018.66'
$ws.Range('E8').Value = '  +2.59%  '
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('D10').Value = "'6.39"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +12.00%  '
$ws.Range('E11').Value = '  +3.84%  '
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('E13').Value = '  +3.80%  '
$ws.Range('D14').Value = "'34.49"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.63%  '
$ws.Range('E15').Value = '  +2.73%  '
$ws.Range('D16').Value = '3.520.10'
$ws.Range('E16').Value = '  +2.34%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '62.461.22'
$ws.Range('E17').Value = '  +1.40%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = "'7.01"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('D19').Value = '3.021.72'
$ws.Range('E19').Value = '  +2.64%  '
$ws.Range('D20').Value = "'447.90"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('E21').Value = '  +2.95%  '
$ws.Range('E22').Value = '  +1.96%  '
$ws.Range('D23').Value = "'7.46"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.82%  '
$ws.Range('D24').Value = "'82.27"
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Value = "'2.23"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.55%  '
$ws.Range('D26').Value = "'10.79"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +14.25%  '
$ws.Range('D27').Value = "'12.03"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('E29').Value = '  +4.11%  '
$ws.Range('D30').Value = "'1.00"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').Value = "'7.18"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.26%  '
$ws.Range('D32').Value = "'2.15"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.79%  '
$ws.Range('D33').Value = "'27.57"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.57%  '
$ws.Range('D34').Value = "'0.109"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.06%  '
$ws.Range('D35').Value = '0.0₃0852'
$ws.Range('E35').Value = '  +10.73%  '
$ws.Range('E36').Value = '  +2.33%  '
$ws.Range('E37').Value = '  +3.66%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').Value = "'3.02"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.79%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = "'2.07"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('D40').Value = "'50.10"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.28%  '
$ws.Range('D41').Value = "'9.02"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.65%  '
$ws.Range('E42').Value = '  +3.24%  '
$ws.Range('E43').Value = '  +8.33%  '
$ws.Range('D44').Value = "'391.91"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.12%  '
$ws.Range('D45').Value = "'40.10"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +9.58%  '
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('D47').Value = '2.751.35'
$ws.Range('E47').Value = '  +1.70%  '
$ws.Range('D48').Value = "'133.33"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.71%  '
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('D51').Value = "'0.107"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.09%  '
